$d = $word.ActiveDocument

# This document contains M2Doc "let" field codes implemented as real Word
# fields (fldChar begin/instrText/separate/.../end). The parser was updated
# to use TokenIteratorFieldRewriterSplit, which stores the field markers as
# plain literal text runs instead: "{" + code, two empty runs (kept so the
# run count/layout matches the previous field's run layout), then "}".
#
# We walk the Fields collection from last to first (so earlier field
# indexes/positions stay valid while we edit), capture each field's code
# text and insertion point, delete the real field, then insert four plain
# runs via InsertXML that reproduce "{<code>", "", "", "}".
#
# The first field (the `let` field) also owns a hidden "_GoBack" bookmark
# that sits right after the field's begin/instrText run in the original
# document; Field.Delete() removes it along with the field runs, so we
# re-create it at the same logical spot (right after the "{<code>" run).

function Escape-Xml($text) {
    $text = $text.Replace("&", "&amp;")
    $text = $text.Replace("<", "&lt;")
    $text = $text.Replace(">", "&gt;")
    return $text
}

$count = $d.Fields.Count

for ($i = $count; $i -ge 1; $i--) {
    $f = $d.Fields($i)
    $code = $f.Code.Text
    $pos = $f.Code.Start - 1
    $needsBookmark = ($i -eq 1)

    $f.Delete()

    $escaped = Escape-Xml($code)

    $xml = '<?xml version="1.0" standalone="yes"?>' +
        '<?mso-application progid="Word.Document"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p>' +
        '<w:r><w:rPr/><w:t>{' + $escaped + '</w:t></w:r>' +
        '<w:r><w:rPr/><w:t/></w:r>' +
        '<w:r><w:rPr/><w:t/></w:r>' +
        '<w:r><w:rPr/><w:t>}</w:t></w:r>' +
        '</w:p></w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

    $ins = $d.Range($pos, $pos)
    $ins.InsertXML($xml)

    if ($needsBookmark) {
        $bmPos = $pos + 1 + $code.Length
        $bmRange = $d.Range($bmPos, $bmPos)
        $d.Bookmarks.Add("_GoBack", $bmRange)
    }
}
